$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalize connector words ("de"/"el"/"del" -> "De"/"El"/"Del") in specific cells ---
$ws.Range("A13").Value = "Ciudad De México"
$ws.Range("A21").Value = "Estado De México"
$ws.Range("B24").Value = "Tlalnepantla De Baz"
$ws.Range("B28").Value = "Apaseo El Grande"
$ws.Range("B33").Value = "Buenavista De Cuéllar"
$ws.Range("B41").Value = "Mineral Del Chico"
$ws.Range("B44").Value = "Atotonilco El Alto"
$ws.Range("B48").Value = "Encarnación De Díaz"
$ws.Range("B55").Value = "Tepatitlán De Morelos"
$ws.Range("B58").Value = "Unión De Tula"
$ws.Range("B60").Value = "Zapotlán El Grande"
$ws.Range("B76").Value = "Ixtlán Del Río"
$ws.Range("B81").Value = "Chalcatongo De Hidalgo"
$ws.Range("B83").Value = "Ixtlán De Juárez"
$ws.Range("B84").Value = "Oaxaca De Juárez"
$ws.Range("B90").Value = "Los Reyes De Juárez"
$ws.Range("B111").Value = "Poza Rica De Hidalgo"

# --- Tiny floating point correction ---
$ws.Range("D64").Value = 0.09142857142857144

# --- Remove trailing metadata/footer rows (117-121 and 476-480) ---
$ws.Rows("476:480").Delete()
$ws.Rows("117:121").Delete()
